$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) stays text so values like "1.000" / "29.401.79" are not
# coerced into numbers/dates by Excel when assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.401.79'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.869.71'
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '243.65'
$ws.Range("D6").Value = '0.7036'
$ws.Range("E6").Value = '  -3.15%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '0.07924'
$ws.Range("E8").Value = '  -0.99%  '
$ws.Range("D9").Value = '0.3133'
$ws.Range("E9").Value = '  -0.90%  '
$ws.Range("D10").Value = '24.48'
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D11").Value = '0.07840'
$ws.Range("E11").Value = '  -4.69%  '
$ws.Range("D12").Value = '1.904.70'
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").Value = '93.85'
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = '5.183'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").Value = '0.7015'
$ws.Range("E15").Value = '  -1.50%  '
$ws.Range("D16").Value = '6.525'
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").Value = '0.000008406'
$ws.Range("E17").Value = '  -0.98%  '
$ws.Range("D18").Value = '29.569.14'
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("D19").Value = '252.65'
$ws.Range("D20").Value = '2.145.25'
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '7.679'
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = '0.1556'
$ws.Range("E25").Value = '  -3.31%  '
$ws.Range("D26").Value = '9.019'
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").Value = '161.46'
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("D28").Value = '18.85'
$ws.Range("E28").Value = '  +1.73%  '
$ws.Range("D29").Value = '1.502'
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").Value = '4.318'
$ws.Range("E30").Value = '  -2.06%  '
$ws.Range("D31").Value = '4.256'
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("E32").Value = '  +2.11%  '
$ws.Range("D33").Value = '0.05266'
$ws.Range("E33").Value = '  -1.81%  '
$ws.Range("D34").Value = '1.899'
$ws.Range("E34").Value = '  -2.15%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.181'
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.7495'
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("D37").Value = '2.711'
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("D38").Value = '0.01881'
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").Value = '1.274.14'
$ws.Range("E39").Value = '  -0.32%  '
$ws.Range("D40").Value = '2.768'
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("D41").Value = '0.8927'
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("D42").Value = '110.17'
$ws.Range("E42").Value = '  -2.84%  '
$ws.Range("D43").Value = '6.037'
$ws.Range("E43").Value = '  -5.97%  '
$ws.Range("D44").Value = '71.08'
$ws.Range("E44").Value = '  -4.26%  '
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("E46").Value = '  -4.94%  '
$ws.Range("D47").Value = '2.030.43'
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("D48").Value = '9.631'
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("D49").Value = '1.801'
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("D50").Value = '0.5182'
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("E51").Value = '  -1.10%  '
